$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells contain price values stored as text (e.g. "34.872.12" uses dots as
# thousands separators, which Excel cannot parse as a single number). Excel's COM
# Value setter auto-converts plain numeric-looking strings into real numbers, so we
# force Text format before assigning, then restore the default "Normal" style so the
# saved cell style matches the original (no explicit style index).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "34.872.12"
$ws.Range("E2").Value = "  -1.28%  "
Set-TextValue $ws.Range("D3") "1.871.34"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").Value = "  -0.94%  "
Set-TextValue $ws.Range("D5") "244.12"
$ws.Range("E5").Value = "  -3.77%  "
Set-TextValue $ws.Range("D6") "0.677"
$ws.Range("E6").Value = "  -6.09%  "
$ws.Range("E7").Value = "  -0.95%  "
Set-TextValue $ws.Range("D8") "42.65"
$ws.Range("E8").Value = "  +4.72%  "
Set-TextValue $ws.Range("D9") "0.342"
$ws.Range("E9").Value = "  -4.93%  "
Set-TextValue $ws.Range("D10") "0.0732"
$ws.Range("E10").Value = "  -2.47%  "
Set-TextValue $ws.Range("D11") "0.0967"
$ws.Range("E11").Value = "  -2.21%  "
Set-TextValue $ws.Range("D12") "12.87"
$ws.Range("E12").Value = "  +1.98%  "
Set-TextValue $ws.Range("D13") "2.142.24"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D15") "4.82"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D16") "1.862.78"
$ws.Range("E16").Value = "  -2.63%  "
Set-TextValue $ws.Range("D17") "34.806.94"
$ws.Range("E17").Value = "  -1.45%  "
Set-TextValue $ws.Range("D18") "72.34"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("E19").Value = "  -4.40%  "
Set-TextValue $ws.Range("D20") "242.86"
$ws.Range("E20").Value = "  -0.40%  "
Set-TextValue $ws.Range("D21") "12.58"
$ws.Range("E21").Value = "  -3.19%  "
Set-TextValue $ws.Range("D22") "4.89"
$ws.Range("E22").Value = "  -3.65%  "
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("E24").Value = "  +5.04%  "
Set-TextValue $ws.Range("D25") "2.15"
$ws.Range("E25").Value = "  -11.57%  "
Set-TextValue $ws.Range("D26") "163.49"
$ws.Range("E26").Value = "  -2.05%  "
Set-TextValue $ws.Range("D27") "8.38"
$ws.Range("E27").Value = "  -2.74%  "
Set-TextValue $ws.Range("D28") "18.03"
$ws.Range("E28").Value = "  -3.69%  "
$ws.Range("E29").Value = "  -4.98%  "
Set-TextValue $ws.Range("D30") "4.128.46"
$ws.Range("E30").Value = "  +0.05%  "
Set-TextValue $ws.Range("D31") "1.75"
$ws.Range("E31").Value = "  +8.17%  "
Set-TextValue $ws.Range("D32") "4.18"
$ws.Range("E32").Value = "  -3.81%  "
Set-TextValue $ws.Range("D33") "0.0571"
$ws.Range("E33").Value = "  -2.04%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D34") "4.13"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D35") "1.00"
$ws.Range("E35").Value = "  -1.00%  "
Set-TextValue $ws.Range("D36") "0.833"
$ws.Range("E36").Value = "  -8.97%  "
$ws.Range("E37").Value = "  -4.45%  "
Set-TextValue $ws.Range("D38") "1.47"
$ws.Range("E38").Value = "  -26.48%  "
Set-TextValue $ws.Range("D39") "97.65"
$ws.Range("E39").Value = "  +0.54%  "
Set-TextValue $ws.Range("D40") "16.97"
$ws.Range("E40").Value = "  -2.06%  "
Set-TextValue $ws.Range("D41") "0.0665"
$ws.Range("E41").Value = "  +3.77%  "
Set-TextValue $ws.Range("D42") "0.0210"
$ws.Range("E42").Value = "  -3.37%  "
$ws.Range("E43").Value = "  -3.76%  "
Set-TextValue $ws.Range("D44") "1.282.36"
$ws.Range("E44").Value = "  -4.13%  "
Set-TextValue $ws.Range("D45") "0.0816"
$ws.Range("E45").Value = "  +10.42%  "
Set-TextValue $ws.Range("D46") "2.31"
$ws.Range("E46").Value = "  -4.94%  "
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("E48").Value = "  -1.45%  "
Set-TextValue $ws.Range("D49") "11.85"
$ws.Range("E49").Value = "  -3.76%  "
Set-TextValue $ws.Range("D50") "6.31"
$ws.Range("E50").Value = "  -6.97%  "
Set-TextValue $ws.Range("D51") "42.62"
$ws.Range("E51").Value = "  -5.47%  "
